$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.234.71"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.855.57"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.08%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7036"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.76"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08013"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +4.56%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3026"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.61"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08207"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.31%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.849.36"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.76%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.196"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.02%  "

$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7060"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.57"
$ws.Range("D15").ClearFormats()

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.125.97"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.822"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.38%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007848"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.62%  "

$ws.Range("E19").Value = "  -0.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.52"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9989"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("B23").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.059.79"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.509"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.86"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.882"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.88%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1415"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.08"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.913"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.399"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.475"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.341"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.021"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05175"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.167"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7132"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9935"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.99%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.673"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01848"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.718"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.35%  "

$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.154.69"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.85%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9329"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.972"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4255"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.61%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.13"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.54"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.63%  "

$ws.Range("E48").Value = "  -3.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.740"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.978.03"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.161"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.15%  "
